$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 0.061724
$ws.Range("H2").Value = 0.185172
$ws.Range("I2").Value = 0.09652262708432048
$ws.Range("J2").Value = 0.09652262708432047
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 0.1445178371693333
$ws.Range("R2").Value = 1.300660534524
$ws.Range("S2").Value = 0.003835339100945612
$ws.Range("T2").Value = 0.003835339100945611
$ws.Range("G3").Value = 0.061724
$ws.Range("H3").Value = 0.185172
$ws.Range("I3").Value = 0.09652262708432048
$ws.Range("J3").Value = 0.09652262708432047
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("Q3").Value = 2.101873401801333
$ws.Range("R3").Value = 18.916860616212
$ws.Range("S3").Value = 0.05578133053375676
$ws.Range("T3").Value = 0.05578133053375675
$ws.Range("G4").Value = 0.061724
$ws.Range("H4").Value = 0.185172
$ws.Range("I4").Value = 0.09652262708432048
$ws.Range("J4").Value = 0.09652262708432047
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("Q4").Value = 1.390638222306667
$ws.Range("R4").Value = 12.51574400076
$ws.Range("S4").Value = 0.03690595744961811
$ws.Range("T4").Value = 0.03690595744961811
$ws.Range("I5").Value = 0.8735221647273214
$ws.Range("J5").Value = 0.8735221647273215
$ws.Range("M5").Value = 2.341355666666667
$ws.Range("N5").Value = 7.024067000000001
$ws.Range("O5").Value = 0.03973512964576821
$ws.Range("P5").Value = 0.0397351296457682
$ws.Range("Q5").Value = 1.307875031784889
$ws.Range("R5").Value = 11.770875286064
$ws.Range("S5").Value = 0.03470951646389221
$ws.Range("T5").Value = 0.03470951646389221
$ws.Range("I6").Value = 0.8735221647273214
$ws.Range("J6").Value = 0.8735221647273215
$ws.Range("O6").Value = 0.5779093692199981
$ws.Range("P6").Value = 0.5779093692199981
$ws.Range("Q6").Value = 19.02178856280355
$ws.Range("S6").Value = 0.5048166432172536
$ws.Range("T6").Value = 0.5048166432172536
$ws.Range("I7").Value = 0.8735221647273214
$ws.Range("J7").Value = 0.8735221647273215
$ws.Range("O7").Value = 0.3823555011342337
$ws.Range("P7").Value = 0.3823555011342337
$ws.Range("S7").Value = 0.3339960050461755
$ws.Range("T7").Value = 0.3339960050461756
$ws.Range("I8").Value = 0.02995520818835809
$ws.Range("J8").Value = 0.02995520818835809
$ws.Range("M8").Value = 2.341355666666667
$ws.Range("N8").Value = 7.024067000000001
$ws.Range("O8").Value = 0.03973512964576821
$ws.Range("P8").Value = 0.0397351296457682
$ws.Range("Q8").Value = 0.04485022869877778
$ws.Range("R8").Value = 0.403652058289
$ws.Range("S8").Value = 0.001190274080930386
$ws.Range("T8").Value = 0.001190274080930386
$ws.Range("I9").Value = 0.02995520818835809
$ws.Range("J9").Value = 0.02995520818835809
$ws.Range("O9").Value = 0.5779093692199981
$ws.Range("P9").Value = 0.5779093692199981
$ws.Range("Q9").Value = 0.652303581434111
$ws.Range("R9").Value = 5.870732232907
$ws.Range("S9").Value = 0.01731139546898775
$ws.Range("T9").Value = 0.01731139546898775
$ws.Range("I10").Value = 0.02995520818835809
$ws.Range("J10").Value = 0.02995520818835809
$ws.Range("O10").Value = 0.3823555011342337
$ws.Range("P10").Value = 0.3823555011342337
$ws.Range("S10").Value = 0.01145353863843996
$ws.Range("T10").Value = 0.01145353863843996
